$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Type" column header in G1, matching the style of the other headers (e.g. E1)
$ws.Range("G1").Value = "Type"
$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add the value for the first data row: "P"
$ws.Range("G2").Value = "P"

# Update the active cell selection to G3, as in the edited workbook
$ws.Range("G3").Select()
